# Fixed possible bug in test level for projectiles.
# Column C (the "requires line of sight" / hit-test flag) for rows 38-89
# on the main data sheet was incorrectly left at 0; it should be 1, same
# as the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the flawed test-data values in column C (rows 38 through 89).
$ws.Range("C38:C89").Value = 1

# Reflect the selection that was active in Excel when the fix was made:
# the corrected range C38:C89, with C38 as the active cell, and the view
# scrolled down so row 58 is at the top of the window.
$ws.Activate()
$ws.Range("C38:C89").Select()
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
